# Updated symbol list on Thu Dec 29 09:25:49 UTC 2022 with GitHub Actions
#
# This script reproduces the published price/volume refresh for the
# cryptos worksheet: most rows just get an updated "Price" (column D)
# figure, while rows 41-44 also got their coin re-ranked (the
# BKEXToken / CEJI / KickToken trio rotated position) which moves the
# Coin name, Link and Volume(1h) text along with the new price.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D ("Price") text updates -------------------------------
# These cells hold their numbers as literal text (t="inlineStr" in the
# original file), so we force a Text number format before writing the
# value to stop Excel from re-interpreting the numeric-looking string
# as a real number (which would also lose trailing/leading zeros).

$priceUpdates = @{
    "D2"  = "245.98"
    "D3"  = "23.93"
    "D4"  = "5.180"
    "D5"  = "0.05733"
    "D6"  = "6.475"
    "D8"  = "0.8133"
    "D9"  = "0.8501"
    "D10" = "0.1370"
    "D11" = "0.06941"
    "D12" = "0.03194"
    "D13" = "0.02879"
    "D14" = "0.09326"
    "D15" = "3.818"
    "D16" = "0.001530"
    "D17" = "0.04685"
    "D18" = "0.0005974"
    "D19" = "0.006211"
    "D20" = "0.001239"
    "D21" = "0.004816"
    "D22" = "0.00008494"
    "D23" = "3.538"
    "D24" = "2.163"
    "D40" = "0.03692"
    "D41" = "0.006386"
    "D42" = "0.1051"
    "D43" = "0.002253"
    "D44" = "0.007807"
    "D45" = "0.00005468"
    "D46" = "0.00000000749"
    "D47" = "0.3997"
    "D48" = "0.002038"
    "D49" = "0.00002098"
    "D50" = "0.0001998"
}

$priceAddrs = @(
    "D2","D3","D4","D5","D6","D8","D9","D10","D11","D12","D13","D14",
    "D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D40",
    "D41","D42","D43","D44","D45","D46","D47","D48","D49","D50"
)

foreach ($addr in $priceAddrs) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    # Flip the format back to General right away (per cell) - doing this
    # as one big batch at the end was observed to be unreliable.
    $cell.NumberFormat = "General"
}

# --- Rows 41-43 re-ranking -------------------------------------------
# BKEXToken (row41) / CEJI (row42) / KickToken (row43) rotated: the coin
# that used to be in row43 moved up to row41, row41's coin moved to
# row42, and row42's coin moved to row43. Column D was already updated
# above; here we fix the Coin, Link and Volume(1h) columns to match.

$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E43").Value = "42CEJICEJI"

# Row 44 (LocalTraders) keeps its Coin/Link, only the "Best in 24h" tag
# on the volume text is dropped now that it's no longer top performer.
$ws.Range("E44").Value = "43LocalTradersLCT"
